$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.805.53"
$ws.Range("E2").Value = "  +4.81%  "

$ws.Range("D3").Value = "1.611.53"
$ws.Range("E3").Value = "  +3.72%  "

$ws.Range("E4").Value = "  -0.50%  "

$ws.Range("D5").Formula = "'213.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.44%  "

$ws.Range("E6").Value = "  +6.99%  "

$ws.Range("D7").Formula = "'0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.49%  "

$ws.Range("D8").Formula = "'26.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +11.32%  "

$ws.Range("E9").Value = "  +3.08%  "

$ws.Range("E10").Value = "  +2.57%  "

$ws.Range("D12").Value = "1.842.11"
$ws.Range("E12").Value = "  +3.75%  "

$ws.Range("D13").Value = "1.632.81"
$ws.Range("E13").Value = "  +5.43%  "

$ws.Range("D14").Value = "29.814.43"
$ws.Range("E14").Value = "  +4.79%  "

$ws.Range("E15").Value = "  +5.41%  "

$ws.Range("E16").Value = "  +3.53%  "

$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").Formula = "'244.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.01%  "

$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Formula = "'63.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.04%  "

$ws.Range("D19").Formula = "'7.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.23%  "

$ws.Range("D20").Value = "0.0₃0696"
$ws.Range("E20").Value = "  +3.42%  "

$ws.Range("D21").Formula = "'0.995"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.44%  "

$ws.Range("E22").Value = "  +4.32%  "

$ws.Range("E23").Value = "  +4.09%  "

$ws.Range("D24").Formula = "'2.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.18%  "

$ws.Range("D25").Formula = "'156.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.45%  "

$ws.Range("E26").Value = "  +4.09%  "

$ws.Range("E27").Value = "  +5.42%  "

$ws.Range("D28").Formula = "'6.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.92%  "

$ws.Range("D29").Formula = "'0.996"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("E30").Value = "  +1.36%  "

$ws.Range("E31").Value = "  +1.06%  "

$ws.Range("E32").Value = "  +3.15%  "

$ws.Range("D33").Value = "1.440.94"
$ws.Range("E33").Value = "  +4.21%  "

$ws.Range("D34").Formula = "'3.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.58%  "

$ws.Range("E35").Value = "  -0.39%  "

$ws.Range("D36").Formula = "'1.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.66%  "

$ws.Range("E37").Value = "  +9.89%  "

$ws.Range("E38").Value = "  +0.55%  "

$ws.Range("E39").Value = "  +3.27%  "

$ws.Range("D40").Formula = "'0.535"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.87%  "

$ws.Range("D41").Formula = "'55.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +28.68%  "

$ws.Range("E42").Value = "  +1.35%  "

$ws.Range("D43").Formula = "'0.799"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.34%  "

$ws.Range("D44").Formula = "'0.995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.41%  "

$ws.Range("D45").Formula = "'0.0467"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.68%  "

$ws.Range("D46").Formula = "'66.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.00%  "

$ws.Range("D47").Formula = "'5.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.88%  "

$ws.Range("D48").Value = "1.752.11"
$ws.Range("E48").Value = "  +3.99%  "

$ws.Range("D49").Formula = "'86.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.08%  "

$ws.Range("E50").Value = "  -4.37%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Formula = "'0.0521"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.95%  "

